$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")

$ws.Cells.Item(26, 8).Value = 0.57278
$ws.Cells.Item(26, 9).Value = 0.02089
$ws.Cells.Item(27, 8).Value = 0.05771
$ws.Cells.Item(27, 9).Value = 0.02182
$ws.Cells.Item(28, 8).Value = 0.5660500000000001
$ws.Cells.Item(28, 9).Value = 0.02004
$ws.Cells.Item(29, 8).Value = 0.03576
$ws.Cells.Item(29, 9).Value = 0.01711
$ws.Cells.Item(30, 8).Value = 0.5732699999999999
$ws.Cells.Item(30, 9).Value = 0.02065
$ws.Cells.Item(31, 8).Value = 0.0567
$ws.Cells.Item(31, 9).Value = 0.02319
$ws.Cells.Item(32, 8).Value = 0.56357
$ws.Cells.Item(32, 9).Value = 0.02053
$ws.Cells.Item(33, 8).Value = 0.04014
$ws.Cells.Item(33, 9).Value = 0.0187
$ws.Cells.Item(34, 8).Value = 0.56101
$ws.Cells.Item(34, 9).Value = 0.01209
$ws.Cells.Item(35, 8).Value = 0.00203
$ws.Cells.Item(35, 9).Value = 0.00433
$ws.Cells.Item(36, 8).Value = 0.56189
$ws.Cells.Item(36, 9).Value = 0.0122
$ws.Cells.Item(37, 8).Value = 0.00203
$ws.Cells.Item(37, 9).Value = 0.00433
$ws.Cells.Item(38, 8).Value = 0.56126
$ws.Cells.Item(38, 9).Value = 0.01213
$ws.Cells.Item(39, 8).Value = 0.00203
$ws.Cells.Item(39, 9).Value = 0.00433
$ws.Cells.Item(40, 8).Value = 0.5620000000000001
$ws.Cells.Item(40, 9).Value = 0.01196
$ws.Cells.Item(41, 8).Value = 0.00203
$ws.Cells.Item(41, 9).Value = 0.00433
$ws.Cells.Item(66, 8).Value = 0.49812
$ws.Cells.Item(66, 9).Value = 0.02334
$ws.Cells.Item(67, 8).Value = 0.01147
$ws.Cells.Item(67, 9).Value = 0.01009
$ws.Cells.Item(68, 8).Value = 0.5091
$ws.Cells.Item(68, 9).Value = 0.02301
$ws.Cells.Item(69, 8).Value = 0.01012
$ws.Cells.Item(69, 9).Value = 0.00984
$ws.Cells.Item(70, 8).Value = 0.50042
$ws.Cells.Item(70, 9).Value = 0.02383
$ws.Cells.Item(71, 8).Value = 0.01114
$ws.Cells.Item(71, 9).Value = 0.01004
$ws.Cells.Item(72, 8).Value = 0.51074
$ws.Cells.Item(72, 9).Value = 0.02438
$ws.Cells.Item(73, 8).Value = 0.01484
$ws.Cells.Item(73, 9).Value = 0.01199
$ws.Cells.Item(74, 8).Value = 0.52986
$ws.Cells.Item(74, 9).Value = 0.01385
$ws.Cells.Item(75, 8).Value = 0.00034
$ws.Cells.Item(75, 9).Value = 0.00166
$ws.Cells.Item(76, 8).Value = 0.53568
$ws.Cells.Item(76, 9).Value = 0.01487
$ws.Cells.Item(77, 8).Value = 0.00168
$ws.Cells.Item(77, 9).Value = 0.00337
$ws.Cells.Item(78, 8).Value = 0.5296999999999999
$ws.Cells.Item(78, 9).Value = 0.01316
$ws.Cells.Item(79, 8).Value = 0.00034
$ws.Cells.Item(79, 9).Value = 0.00166
$ws.Cells.Item(80, 8).Value = 0.53483
$ws.Cells.Item(80, 9).Value = 0.01505
$ws.Cells.Item(81, 8).Value = 0.00202
$ws.Cells.Item(81, 9).Value = 0.0036
$ws.Cells.Item(106, 8).Value = 0.57143
$ws.Cells.Item(106, 9).Value = 0.02359
$ws.Cells.Item(107, 8).Value = 0.05533
$ws.Cells.Item(107, 9).Value = 0.02438
$ws.Cells.Item(108, 8).Value = 0.56487
$ws.Cells.Item(108, 9).Value = 0.02316
$ws.Cells.Item(109, 8).Value = 0.03141
$ws.Cells.Item(109, 9).Value = 0.0235
$ws.Cells.Item(110, 8).Value = 0.57156
$ws.Cells.Item(110, 9).Value = 0.02389
$ws.Cells.Item(111, 8).Value = 0.05567
$ws.Cells.Item(111, 9).Value = 0.02505
$ws.Cells.Item(112, 8).Value = 0.5628
$ws.Cells.Item(112, 9).Value = 0.02273
$ws.Cells.Item(113, 8).Value = 0.03478
$ws.Cells.Item(113, 9).Value = 0.02348
$ws.Cells.Item(114, 8).Value = 0.56606
$ws.Cells.Item(114, 9).Value = 0.01672
$ws.Cells.Item(115, 8).Value = 0.00135
$ws.Cells.Item(115, 9).Value = 0.0031
$ws.Cells.Item(116, 8).Value = 0.56667
$ws.Cells.Item(116, 9).Value = 0.01658
$ws.Cells.Item(117, 8).Value = 0.00203
$ws.Cells.Item(117, 9).Value = 0.0036
$ws.Cells.Item(118, 8).Value = 0.56594
$ws.Cells.Item(118, 9).Value = 0.01678
$ws.Cells.Item(119, 8).Value = 0.00135
$ws.Cells.Item(119, 9).Value = 0.0031
$ws.Cells.Item(120, 8).Value = 0.56649
$ws.Cells.Item(120, 9).Value = 0.01683
$ws.Cells.Item(121, 8).Value = 0.00203
$ws.Cells.Item(121, 9).Value = 0.0036
$ws.Cells.Item(146, 8).Value = 0.49723
$ws.Cells.Item(146, 9).Value = 0.02511
$ws.Cells.Item(147, 8).Value = 0.01417
$ws.Cells.Item(147, 9).Value = 0.01138
$ws.Cells.Item(148, 8).Value = 0.50467
$ws.Cells.Item(148, 9).Value = 0.0223
$ws.Cells.Item(149, 8).Value = 0.01114
$ws.Cells.Item(149, 9).Value = 0.01112
$ws.Cells.Item(150, 8).Value = 0.4998
$ws.Cells.Item(150, 9).Value = 0.02466
$ws.Cells.Item(151, 8).Value = 0.01485
$ws.Cells.Item(151, 9).Value = 0.01269
$ws.Cells.Item(152, 8).Value = 0.50775
$ws.Cells.Item(152, 9).Value = 0.0237
$ws.Cells.Item(153, 8).Value = 0.01788
$ws.Cells.Item(153, 9).Value = 0.01588
$ws.Cells.Item(154, 8).Value = 0.53355
$ws.Cells.Item(154, 9).Value = 0.01769
$ws.Cells.Item(155, 8).Value = 0.00135
$ws.Cells.Item(155, 9).Value = 0.00309
$ws.Cells.Item(156, 8).Value = 0.53928
$ws.Cells.Item(156, 9).Value = 0.01764
$ws.Cells.Item(157, 8).Value = 0.00473
$ws.Cells.Item(157, 9).Value = 0.00927
$ws.Cells.Item(158, 8).Value = 0.53339
$ws.Cells.Item(158, 9).Value = 0.01714
$ws.Cells.Item(159, 8).Value = 0.00135
$ws.Cells.Item(159, 9).Value = 0.00309
$ws.Cells.Item(160, 8).Value = 0.53854
$ws.Cells.Item(160, 9).Value = 0.01757
$ws.Cells.Item(161, 8).Value = 0.00473
$ws.Cells.Item(161, 9).Value = 0.00927
